$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 42.0625
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 137.66667
$ws.Range("K2").Value = 20
$ws.Range("L2").Value = 137.66667
$ws.Range("M2").Value = 93
$ws.Range("N2").Value = -363.66667

$ws.Range("H12").Value = 4001
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H28").Value = 1172
$ws.Range("I28").Value = 946.38464
$ws.Range("J28").Value = 2149.6667
$ws.Range("K28").Value = 946.38464
$ws.Range("L28").Value = 2149.6667
$ws.Range("M28").Value = -461.38464
$ws.Range("N28").Value = -3119.6667

$ws.Range("H116").Value = 25006780
$ws.Range("I116").Value = 41670096
$ws.Range("J116").Value = 11801.5
$ws.Range("K116").Value = 41670096
$ws.Range("L116").Value = 11801.5
$ws.Range("M116").Value = -41666654
$ws.Range("N116").Value = -18685.5

$ws.Range("H138").Value = 3924.2415
$ws.Range("I138").Value = 1122.5
$ws.Range("J138").Value = 5253.8813
$ws.Range("K138").Value = 3367.5
$ws.Range("L138").Value = 15761.6439
$ws.Range("M138").Value = 1772.5
$ws.Range("N138").Value = -26041.6439

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4050.9092
$ws.Range("I2").Value = 1234.8334
$ws.Range("J2").Value = 7430.2
$ws.Range("K2").Value = 1234.8334
$ws.Range("L2").Value = 7430.2
$ws.Range("M2").Value = -1121.8334
$ws.Range("N2").Value = -7656.2

$ws.Range("H45").Value = 11500.5
$ws.Range("J45").Value = 13000
$ws.Range("L45").Value = 13000
$ws.Range("N45").Value = -13754

$ws.Range("H61").Value = 5139.9346
$ws.Range("I61").Value = 2898.4082
$ws.Range("K61").Value = 2898.4082
$ws.Range("M61").Value = -2686.4082

$ws.Range("H63").Value = 2996
$ws.Range("I63").Value = 2996
$ws.Range("K63").Value = 2996
$ws.Range("M63").Value = -2310

$ws.Range("H66").Value = 2996
$ws.Range("I66").Value = 2996
$ws.Range("K66").Value = 14980
$ws.Range("M66").Value = -11548

$ws.Range("H102").Value = 920.82355
$ws.Range("I102").Value = 914.75
$ws.Range("K102").Value = 914.75
$ws.Range("M102").Value = 707.25

$ws.Range("H116").Value = 4050.9092
$ws.Range("I116").Value = 1234.8334
$ws.Range("J116").Value = 7430.2
$ws.Range("K116").Value = 1234.8334
$ws.Range("L116").Value = 7430.2
$ws.Range("M116").Value = 1059.1666
$ws.Range("N116").Value = -12018.2

$ws.Range("H132").Value = 4732.184
$ws.Range("I132").Value = 1330.05
$ws.Range("K132").Value = 3990.15
$ws.Range("M132").Value = -1460.15

$ws.Range("H136").Value = 5139.9346
$ws.Range("I136").Value = 2898.4082
$ws.Range("K136").Value = 8695.2246
$ws.Range("M136").Value = -6145.2246

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4050.9092
$ws.Range("I3").Value = 1234.8334
$ws.Range("J3").Value = 7430.2
$ws.Range("K3").Value = 1234.8334
$ws.Range("L3").Value = 7430.2
$ws.Range("M3").Value = -1120.8334
$ws.Range("N3").Value = -7658.2

$ws.Range("H107").Value = 102277390
$ws.Range("I107").Value = 281250700
$ws.Range("K107").Value = 281250700
$ws.Range("M107").Value = -281248780

$ws.Range("H134").Value = 5704.021
$ws.Range("I134").Value = 1994.0416
$ws.Range("K134").Value = 5982.1248
$ws.Range("M134").Value = -3447.1248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6043.1763
$ws.Range("I16").Value = 4418.1665
$ws.Range("J16").Value = 6929.5454
$ws.Range("K16").Value = 4418.1665
$ws.Range("L16").Value = 6929.5454
$ws.Range("M16").Value = -4131.1665
$ws.Range("N16").Value = -7503.5454

$ws.Range("H22").Value = 127.4
$ws.Range("I22").Value = 84.5
$ws.Range("J22").Value = 299
$ws.Range("K22").Value = 84.5
$ws.Range("L22").Value = 299
$ws.Range("M22").Value = 265.5
$ws.Range("N22").Value = -999

$ws.Range("H31").Value = 5562.6387
$ws.Range("I31").Value = 2305.761
$ws.Range("J31").Value = 11324.808
$ws.Range("K31").Value = 2305.761
$ws.Range("L31").Value = 11324.808
$ws.Range("M31").Value = -2010.761
$ws.Range("N31").Value = -11914.808

$ws.Range("H34").Value = 5562.6387
$ws.Range("I34").Value = 2305.761
$ws.Range("J34").Value = 11324.808
$ws.Range("K34").Value = 2305.761
$ws.Range("L34").Value = 11324.808
$ws.Range("M34").Value = -2103.761
$ws.Range("N34").Value = -11728.808

$ws.Range("H58").Value = 6853138
$ws.Range("I58").Value = 9616654
$ws.Range("J58").Value = 10144.619
$ws.Range("K58").Value = 9616654
$ws.Range("L58").Value = 10144.619
$ws.Range("M58").Value = -9616451
$ws.Range("N58").Value = -10550.619

$ws.Range("H86").Value = 10424500
$ws.Range("I86").Value = 20839002
$ws.Range("J86").Value = 9998.666999999999
$ws.Range("K86").Value = 20839002
$ws.Range("L86").Value = 9998.666999999999
$ws.Range("M86").Value = -20837879
$ws.Range("N86").Value = -12244.667

$ws.Range("H89").Value = 10424500
$ws.Range("I89").Value = 20839002
$ws.Range("J89").Value = 9998.666999999999
$ws.Range("K89").Value = 104195010
$ws.Range("L89").Value = 49993.335
$ws.Range("M89").Value = -104189394
$ws.Range("N89").Value = -61225.335

$ws.Range("H93").Value = 16874.834
$ws.Range("I93").Value = 9461.4
$ws.Range("J93").Value = 53942
$ws.Range("K93").Value = 9461.4
$ws.Range("L93").Value = 53942
$ws.Range("M93").Value = -7589.4
$ws.Range("N93").Value = -57686

$ws.Range("H99").Value = 7661.2
$ws.Range("I99").Value = 3799
$ws.Range("K99").Value = 3799
$ws.Range("M99").Value = -2301

$ws.Range("H105").Value = 4468593
$ws.Range("I105").Value = 6494907.5
$ws.Range("K105").Value = 6494907.5
$ws.Range("M105").Value = -6493160.5

$ws.Range("H107").Value = 3797.2
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 3797.2
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 3797.2
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -7637.2

$ws.Range("H113").Value = 6043.1763
$ws.Range("I113").Value = 4418.1665
$ws.Range("J113").Value = 6929.5454
$ws.Range("K113").Value = 4418.1665
$ws.Range("L113").Value = 6929.5454
$ws.Range("M113").Value = -2248.1665
$ws.Range("N113").Value = -11269.5454

$ws.Range("H126").Value = 7661.2
$ws.Range("I126").Value = 3799
$ws.Range("K126").Value = 11397
$ws.Range("M126").Value = -8927

$ws.Range("H132").Value = 3736.884
$ws.Range("I132").Value = 1550.5714
$ws.Range("J132").Value = 9093.35
$ws.Range("K132").Value = 4651.7142
$ws.Range("L132").Value = 27280.05
$ws.Range("M132").Value = -2121.7142
$ws.Range("N132").Value = -32340.05

$ws.Range("H134").Value = 3703.111
$ws.Range("I134").Value = 1530.862
$ws.Range("J134").Value = 9180.956
$ws.Range("K134").Value = 4592.586
$ws.Range("L134").Value = 27542.868
$ws.Range("M134").Value = -2057.586
$ws.Range("N134").Value = -32612.868

$ws.Range("H136").Value = 6853138
$ws.Range("I136").Value = 9616654
$ws.Range("J136").Value = 10144.619
$ws.Range("K136").Value = 28849962
$ws.Range("L136").Value = 30433.857
$ws.Range("M136").Value = -28847412
$ws.Range("N136").Value = -35533.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1757.9
$ws.Range("I131").Value = 1597.3334
$ws.Range("K131").Value = 4792.0002
$ws.Range("M131").Value = 247.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2522.75
$ws.Range("I31").Value = 191
$ws.Range("K31").Value = 191
$ws.Range("M31").Value = 101

$ws.Range("H35").Value = 26249.25
$ws.Range("J35").Value = 29999
$ws.Range("L35").Value = 29999
$ws.Range("N35").Value = -30595

$ws.Range("H37").Value = 2522.75
$ws.Range("I37").Value = 191
$ws.Range("K37").Value = 191
$ws.Range("M37").Value = 86

$ws.Range("H113").Value = 7045.343
$ws.Range("I113").Value = 3354.4546
$ws.Range("K113").Value = 3354.4546
$ws.Range("M113").Value = -1184.4546

$ws.Range("H122").Value = 2654577
$ws.Range("I122").Value = 3256844.5
$ws.Range("K122").Value = 9770533.5
$ws.Range("M122").Value = -9768083.5

$ws.Range("H126").Value = 9995.75
$ws.Range("J126").Value = 9998.666999999999
$ws.Range("L126").Value = 29996.001
$ws.Range("N126").Value = -34936.001

$ws.Range("H132").Value = 8992.5625
$ws.Range("I132").Value = 2839
$ws.Range("K132").Value = 8517
$ws.Range("M132").Value = -5987

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 18521800
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 60190.332
$ws.Range("I62").Value = 171331.33
$ws.Range("J62").Value = 4619.8335
$ws.Range("K62").Value = 171331.33
$ws.Range("L62").Value = 4619.8335
$ws.Range("M62").Value = -170707.33
$ws.Range("N62").Value = -5867.8335

$ws.Range("H65").Value = 60190.332
$ws.Range("I65").Value = 171331.33
$ws.Range("J65").Value = 4619.8335
$ws.Range("K65").Value = 856656.6499999999
$ws.Range("L65").Value = 23099.1675
$ws.Range("M65").Value = -853536.6499999999
$ws.Range("N65").Value = -29339.1675

$ws.Range("H122").Value = 2924.0264
$ws.Range("I122").Value = 1416.4348
$ws.Range("K122").Value = 4249.3044
$ws.Range("M122").Value = -1799.3044

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 9442296
$ws.Range("I132").Value = 12198798
$ws.Range("K132").Value = 36596394
$ws.Range("M132").Value = -36596394
